# The commit swaps the OOXML content of ppt/theme/theme1.xml (the
# presentation's "Integral" theme, used by the slide master / all the
# visible slides) and ppt/theme/theme2.xml (the "Office Theme", used by
# the notes master): theme1.xml ends up holding the Office color scheme
# and theme2.xml ends up holding the Integral color scheme.
#
# Through the PowerPoint object model, a deck's colour scheme is reached
# via <Master/Slide/...>.Theme.ThemeColorScheme, whose 12 slots map 1:1
# onto clrScheme's dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink (in that
# order), each exposed as a ThemeColor with a settable .RGB (a standard
# COM RGB value, i.e. 0xBBGGRR).
#
# Apply the "Office Theme" colour values (previously only present in
# theme2.xml) onto the presentation's theme colour scheme.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $b * 65536 + $g * 256 + $r
}

# clrScheme slot order exposed by ThemeColorScheme.Colors(1..12)
$officeThemeColors = @(
    (RGBVal 0x00 0x00 0x00), # 1  dk1
    (RGBVal 0xFF 0xFF 0xFF), # 2  lt1
    (RGBVal 0x44 0x54 0x6A), # 3  dk2
    (RGBVal 0xE7 0xE6 0xE6), # 4  lt2
    (RGBVal 0x5B 0x9B 0xD5), # 5  accent1
    (RGBVal 0xED 0x7D 0x31), # 6  accent2
    (RGBVal 0xA5 0xA5 0xA5), # 7  accent3
    (RGBVal 0xFF 0xC0 0x00), # 8  accent4
    (RGBVal 0x44 0x72 0xC4), # 9  accent5
    (RGBVal 0x70 0xAD 0x47), # 10 accent6
    (RGBVal 0x05 0x63 0xC1), # 11 hlink
    (RGBVal 0x95 0x4F 0x72)  # 12 folHlink
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
